$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column (D) whose value is being updated below. Force them to
# text ("@") format first so numeric-looking strings (e.g. "1.00", "20.90") keep
# their literal text instead of being normalized into a number by Excel.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D12", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the refreshed coin data (prices / 1h volume %, plus the Aave/Maker row swap).
$ws.Range('D2').Value = '58.949.78'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '2.597.61'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '554.77'
$ws.Range('E5').Value = '  +3.20%  '
$ws.Range('D6').Value = '144.28'
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.602'
$ws.Range('E8').Value = '  +5.12%  '
$ws.Range('D9').Value = '6.79'
$ws.Range('E9').Value = '  +1.83%  '
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('E11').Value = '  +4.87%  '
$ws.Range('D12').Value = '0.337'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('D13').Value = '3.051.05'
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').Value = '58.878.97'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').Value = '20.90'
$ws.Range('E15').Value = '  -1.51%  '
$ws.Range('D16').Value = '2.582.96'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('D18').Value = '4.45'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').Value = '10.08'
$ws.Range('E20').Value = '  -3.16%  '
$ws.Range('D21').Value = '6.16'
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '66.39'
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('D24').Value = '0.428'
$ws.Range('E24').Value = '  +2.88%  '
$ws.Range('D25').Value = '0.997'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = '0.159'
$ws.Range('E26').Value = '  -2.87%  '
$ws.Range('D27').Value = '7.14'
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('D28').Value = '0.0₃0756'
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = '1.68'
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('D31').Value = '5.95'
$ws.Range('E31').Value = '  +1.67%  '
$ws.Range('D32').Value = '153.09'
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('D33').Value = '18.94'
$ws.Range('E33').Value = '  +0.22%  '
$ws.Range('E34').Value = '  -1.97%  '
$ws.Range('D35').Value = '0.874'
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('D36').Value = '1.12'
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').Value = '36.99'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').Value = '1.46'
$ws.Range('E38').Value = '  +1.44%  '
$ws.Range('D39').Value = '0.825'
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('D40').Value = '3.61'
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('D41').Value = '283.81'
$ws.Range('E41').Value = '  -1.66%  '
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').Value = '0.600'
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('D44').Value = '0.0956'
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('D45').Value = '10.64'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').Value = '0.0533'
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('D47').Value = '0.0227'
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '118.78'
$ws.Range('E48').Value = '  +7.25%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '1.917.56'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('D50').Value = '4.44'
$ws.Range('E50').Value = '  -2.62%  '
$ws.Range('D51').Value = '17.84'
$ws.Range('E51').Value = '  -2.56%  '
